# Added Support for Vertical Certificate & Added GUI
# - Fix the certificate number sample value/format
# - Lock the (now dynamic) certificate rows so the template can't be
#   accidentally edited by the GUI-driven generator
# - Leave the cursor parked on the first data row used by the new flow

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the sample certificate number's format
$ws.Range("B2").Value = "EYEQ123/24/25"

# The rows below the header are now filled dynamically by the GUI; give
# them their own explicit (locked) formatting instead of relying on the
# sheet's default style.
$ws.Range("A3:B7").VerticalAlignment = -4160
$ws.Range("A3:B7").Locked = $true

# Park the selection on the newly-protected data range
$ws.Range("B6").Select()
